$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vmPuValues = @{
    "B2" = 1.02; "C2" = 1.036694630684902; "D2" = 1.044345466682616; "E2" = 1.040343640754303; "F2" = 1.05287886545348; "I2" = 1.036917892748111; "J2" = 1.041801549840452; "K2" = 1.047116777370738; "L2" = 1.043126257244017; "M2" = 1.055626379067381; "N2" = 1.043281027537533
    "B3" = 1.02; "C3" = 1.037795243449714; "D3" = 1.045196586551238; "E3" = 1.041392913115981; "F3" = 1.05386332215236; "I3" = 1.037128078460867; "J3" = 1.042545472805016; "K3" = 1.047779012258042; "L3" = 1.043985305131532; "M3" = 1.056423329363004; "N3" = 1.044026006958036
    "B4" = 1.02; "C4" = 1.038507660308317; "D4" = 1.04574733384831; "E4" = 1.04207246211418; "F4" = 1.05450065787718; "I4" = 1.037262759522156; "J4" = 1.04302653782935; "K4" = 1.048206895907588; "L4" = 1.044541166555784; "M4" = 1.05693871936551; "N4" = 1.044507755149884
    "B5" = 1.02; "C5" = 1.038807219804913; "D5" = 1.045978871294667; "E5" = 1.042358287906623; "F5" = 1.054768671517106; "I5" = 1.037319062454998; "J5" = 1.043228704920893; "K5" = 1.048386627874178; "L5" = 1.0447748505561; "M5" = 1.057155319506182; "N5" = 1.044710209341886
    "B6" = 1.02; "C6" = 1.03885752065934; "D6" = 1.046017747628492; "E6" = 1.042406287719062; "F6" = 1.054813676730168; "I6" = 1.037328497370603; "J6" = 1.043262645401337; "K6" = 1.04841679683844; "L6" = 1.044814087115579; "M6" = 1.057191683518413; "N6" = 1.044744198021704
    "B7" = 1.02; "C7" = 1.038511662802909; "D7" = 1.045750427650739; "E7" = 1.042076280769159; "F7" = 1.054504238783795; "I7" = 1.03726351309051; "J7" = 1.043029239481534; "K7" = 1.048209298085965; "L7" = 1.044544289053568; "M7" = 1.056941613861972; "N7" = 1.044510460638723
    "B8" = 1.02; "C8" = 1.037066537037434; "D8" = 1.044633103045186; "E8" = 1.040698123012841; "F8" = 1.053211499352435; "I8" = 1.036989199827921; "J8" = 1.042053024721125; "K8" = 1.047340712308396; "L8" = 1.043416576414529; "M8" = 1.055895772008655; "N8" = 1.043532859541382
    "B9" = 1.02; "C9" = 1.034521914437976; "D9" = 1.042664377570505; "E9" = 1.038274231779615; "F9" = 1.050936046368666; "I9" = 1.036495700674771; "J9" = 1.040330495491736; "K9" = 1.045805367461924; "L9" = 1.041429408525136; "M9" = 1.054050658472512; "N9" = 1.04180788412291
    "B10" = 1.02; "C10" = 1.032826733179371; "D10" = 1.041352017735866; "E10" = 1.036661402194636; "F10" = 1.049420794733867; "I10" = 1.036159906395488; "J10" = 1.039180592078599; "K10" = 1.044778602536688; "L10" = 1.040104639320167; "M10" = 1.052819117862176; "N10" = 1.040656347715031
    "B11" = 1.02; "C11" = 1.03209298698684; "D11" = 1.040783785906197; "E11" = 1.035963764298713; "F11" = 1.048765085425995; "I11" = 1.03601289435011; "J11" = 1.038682303791616; "K11" = 1.044333245297995; "L11" = 1.039531002191028; "M11" = 1.052285502928805; "N11" = 1.040157351801525
    "B12" = 1.02; "C12" = 1.031820482142831; "D12" = 1.040572723932982; "E12" = 1.035704739811076; "F12" = 1.04852158682569; "I12" = 1.035958045637251; "J12" = 1.038497161303543; "K12" = 1.044167705552284; "L12" = 1.039317927248497; "M12" = 1.052087242208168; "N12" = 1.039971946389883
    "B13" = 1.02; "C13" = 1.03187893350507; "D13" = 1.040617997228011; "E13" = 1.035760296498629; "F13" = 1.048573815336378; "I13" = 1.035969821816489; "J13" = 1.038536877532315; "K13" = 1.04420321955706; "L13" = 1.039363632573991; "M13" = 1.052129772173167; "N13" = 1.040011719020256
    "B14" = 1.02; "C14" = 1.032070460815631; "D14" = 1.040766339361338; "E14" = 1.03594235102474; "F14" = 1.048744956509319; "I14" = 1.03600836547221; "J14" = 1.038667000999203; "K14" = 1.044319564060626; "L14" = 1.039513389361884; "M14" = 1.052269115703261; "N14" = 1.040142027277392
    "B15" = 1.02; "C15" = 1.032188472547065; "D15" = 1.040857738468703; "E15" = 1.036054535306412; "F15" = 1.048850410328814; "I15" = 1.036032081441184; "J15" = 1.038747166935064; "K15" = 1.044391232615462; "L15" = 1.039605659383618; "M15" = 1.052354962915844; "N15" = 1.040222307058077
    "B16" = 1.02; "C16" = 1.032875435307045; "D16" = 1.041389730042224; "E16" = 1.036707717454699; "F16" = 1.049464320546922; "I16" = 1.036169629171635; "J16" = 1.039213653976978; "K16" = 1.04480814340242; "L16" = 1.040142709651; "M16" = 1.052854524773857; "N16" = 1.040689456565097
    "B17" = 1.02; "C17" = 1.033306422975452; "D17" = 1.041723442199417; "E17" = 1.03711763619901; "F17" = 1.049849518711804; "I17" = 1.036255478032978; "J17" = 1.039506169193563; "K17" = 1.045069456809028; "L17" = 1.040479586154458; "M17" = 1.053167793216315; "N17" = 1.040982387186841
    "B18" = 1.02; "C18" = 1.033557837705063; "D18" = 1.041918093596177; "E18" = 1.037356805109817; "F18" = 1.050074237311162; "I18" = 1.03630539678366; "J18" = 1.039676752431878; "K18" = 1.045221803039569; "L18" = 1.040676080203793; "M18" = 1.053350483644421; "N18" = 1.041153212672921
    "B19" = 1.02; "C19" = 1.033643568216486; "D19" = 1.041984465147284; "E19" = 1.037438367444463; "F19" = 1.050150867158667; "I19" = 1.036322391418847; "J19" = 1.039734910814062; "K19" = 1.045273736703374; "L19" = 1.040743079545438; "M19" = 1.053412770628878; "N19" = 1.041211453646679
    "B20" = 1.02; "C20" = 1.033260179298313; "D20" = 1.041687637757163; "E20" = 1.037073648532109; "F20" = 1.04980818655564; "I20" = 1.036246283339317; "J20" = 1.039474788796823; "K20" = 1.045041427950661; "L20" = 1.04044344252687; "M20" = 1.053134185967248; "N20" = 1.040950962226339
    "B21" = 1.02; "C21" = 1.032014059665505; "D21" = 1.04072265617001; "E21" = 1.035888737488093; "F21" = 1.048694558012772; "I21" = 1.035997022002695; "J21" = 1.03862868442872; "K21" = 1.044285306656599; "L21" = 1.039469289736794; "M21" = 1.052228083951166; "N21" = 1.040103656292984
    "B22" = 1.02; "C22" = 1.031230812612329; "D22" = 1.040115960888184; "E22" = 1.035144369161927; "F22" = 1.047994728778994; "I22" = 1.035838902049687; "J22" = 1.038096380815483; "K22" = 1.043809242309779; "L22" = 1.038856797871697; "M22" = 1.051658078706712; "N22" = 1.03957059674756
    "B23" = 1.02; "C23" = 1.031646004377262; "D23" = 1.040437578873684; "E23" = 1.035538912988864; "F23" = 1.048365687935326; "I23" = 1.035922857038254; "J23" = 1.038378595785709; "K23" = 1.044061675627189; "L23" = 1.03918149173888; "M23" = 1.051960277931565; "N23" = 1.039853212495414
    "B24" = 1.02; "C24" = 1.033281074742047; "D24" = 1.041703816234892; "E24" = 1.037093524448165; "F24" = 1.049826862657206; "I24" = 1.036250438505842; "J24" = 1.03948896835827; "K24" = 1.04505409321215; "L24" = 1.040459774277227; "M24" = 1.053149371740691; "N24" = 1.04096516192439
    "B25" = 1.02; "C25" = 1.035179541019561; "D25" = 1.043173320481885; "E25" = 1.038900320172997; "F25" = 1.051524004577992; "I25" = 1.036624480557634; "J25" = 1.040776084276339; "K25" = 1.046202856394679; "L25" = 1.041943137802525; "M25" = 1.05452792411519; "N25" = 1.042254105694697
}

foreach ($cellRef in $vmPuValues.Keys) {
    $ws.Range($cellRef).Value = $vmPuValues[$cellRef]
}
